$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post about "木に登れば熊の危険を避けられるとお考えなら..." (row 775) was
# removed from the source data. Delete that entire row and shift everything
# below it up by one (Excel automatically renumbers/shifts dependent rows).
$ws.Rows.Item(775).Delete()
